$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replacement financial figures for rows 2-6 (source data was re-pulled; see
# commit "error solve ifrs list"). Keys are column letters, values are the
# corrected numbers for that row.
$newRowValues = @{
    2 = @{ "D" = 6163; "E" = 138; "F" = 120; "G" = -33; "H" = -44; "I" = -29; "J" = -15; "K" = 6518; "L" = 5191; "M" = 1327; "N" = 992; "O" = 335; "P" = 242; "Q" = 270; "R" = -95; "S" = -29; "T" = 328; "U" = -58; "V" = 3831; "W" = 2.24; "X" = -0.71; "Y" = -2.8; "Z" = -0.67; "AA" = 391.29; "AB" = 306.05; "AC" = -102; "AD" = -16.33; "AE" = 3582; "AF" = 0.47; "AG" = 9; "AH" = 0.51; "AI" = -21.15; "AJ" = 25218517 }
    3 = @{ "D" = 6298; "E" = 246; "F" = 246; "G" = -53; "H" = -172; "I" = -137; "J" = -35; "K" = 6309; "L" = 5148; "M" = 1160; "N" = 818; "O" = 342; "P" = 268; "Q" = 715; "R" = -365; "S" = -538; "T" = 412; "U" = 304; "V" = 3580; "W" = 3.91; "X" = -2.73; "Y" = -15.09; "Z" = -2.68; "AA" = 443.7; "AB" = 203.05; "AC" = -492; "AD" = -4.56; "AE" = 2956; "AF" = 0.76; "AG" = 9; "AH" = 0.38; "AI" = -1.41; "AJ" = 25218517 }
    4 = @{ "D" = 5626; "E" = 257; "F" = 203; "G" = 126; "H" = 146; "I" = 144; "J" = 2; "K" = 5885; "L" = 4394; "M" = 1491; "N" = 1151; "O" = 340; "P" = 378; "Q" = 272; "R" = 235; "S" = -370; "T" = 317; "U" = -45; "V" = 3055; "W" = 4.58; "X" = 2.59; "Y" = 14.64; "Z" = 2.39; "AA" = 294.77; "AB" = 201.58; "AC" = 459; "AD" = 4.45; "AE" = 3058; "AF" = 0.67; "AG" = 10; "AH" = 0.49; "AI" = 3.07; "AJ" = 32523970 }
    5 = @{ "D" = 4670; "E" = 146; "F" = 146; "G" = -159; "H" = -94; "I" = -101; "J" = 8; "K" = 5425; "L" = 4158; "M" = 1267; "N" = 1036; "O" = 231; "P" = 378; "Q" = 193; "R" = -167; "S" = -96; "T" = 98; "U" = 95; "V" = 3048; "W" = 3.14; "X" = -2; "Y" = -9.27; "Z" = -1.65; "AA" = 328.28; "AB" = 184.81; "AC" = -268; "AD" = -5.35; "AE" = 2751; "AF" = 0.52; "AG" = 10; "AH" = 0.7; "AI" = -3.43; "AJ" = 32523970 }
    6 = @{ "D" = 5264; "E" = 142; "F" = 142; "G" = 171; "H" = 158; "I" = 153; "K" = 5541; "L" = 4239; "M" = 1302; "N" = 1050; "P" = 381; "Q" = 285; "R" = -217; "S" = -147; "T" = 132; "U" = 153; "V" = 2697; "W" = 2.69; "X" = 3; "Y" = 14.67; "Z" = 2.88; "AA" = 325.61; "AB" = 206.49; "AC" = 444; "AD" = 3.09; "AE" = 3154; "AF" = 0.43; "AG" = 10; "AH" = 0.73; "AI" = 1.54; "AJ" = 33465588 }
}

foreach ($rowNum in $newRowValues.Keys) {
    $rowCols = $newRowValues[$rowNum]
    foreach ($col in $rowCols.Keys) {
        $ws.Range("$col$rowNum").Value = $rowCols[$col]
    }
}

# Rows 7-9 no longer carry any figures beyond the label columns (A-C); clear
# the previously populated cells so they are removed entirely.
$ws.Range("D7:AI9").ClearContents()
